# Updated cryptos list on Sun Jun 16 17:49:41 UTC 2024 with GitHub Actions
#
# Applies the latest price / 1h-volume-change snapshot to the cryptos
# worksheet. Column D (Price) and column E (Volume(1h)) are stored as
# plain text in the workbook (values like "66.570.22" use '.' as a
# thousands separator, and the percentages keep their padding spaces),
# so every write below is apostrophe-prefixed to force Excel to keep it
# as text instead of re-interpreting it as a number/date, and the style
# is reset to "Normal" afterwards so the quote-prefix flag doesn't stick
# around as a new cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, new price (D), new 1h change (E) -- only the keys that changed
# are present for each row.
$rows = @(
    @{ Row = 2;  D = "66.570.22";   E = "  +0.41%  " },
    @{ Row = 3;  D = "3.593.90";    E = "  +0.69%  " },
    @{ Row = 4;  E = "  +0.09%  " },
    @{ Row = 5;  D = "609.02";      E = "  +0.18%  " },
    @{ Row = 6;  D = "148.20";      E = "  +1.89%  " },
    @{ Row = 7;  E = "  -0.01%  " },
    @{ Row = 8;  E = "  +0.58%  " },
    @{ Row = 9;  E = "  -0.30%  " },
    @{ Row = 10; D = "8.02";        E = "  -0.02%  " },
    @{ Row = 11; E = "  +0.68%  " },
    @{ Row = 12; D = "4.210.95";    E = "  +0.83%  " },
    @{ Row = 13; D = "0.0000209";   E = "  +0.28%  " },
    @{ Row = 14; D = "30.00";       E = "  -0.71%  " },
    @{ Row = 15; D = "3.608.27";    E = "  +1.15%  " },
    @{ Row = 16; D = "66.659.01";   E = "  +0.42%  " },
    @{ Row = 17; E = "  +0.88%  " },
    @{ Row = 18; D = "11.55";       E = "  +0.35%  " },
    @{ Row = 19; D = "6.35";        E = "  +2.20%  " },
    @{ Row = 20; D = "15.06";       E = "  +0.59%  " },
    @{ Row = 21; D = "428.49";      E = "  -0.82%  " },
    @{ Row = 22; D = "0.620";       E = "  +1.78%  " },
    @{ Row = 23; D = "79.07";       E = "  +0.38%  " },
    @{ Row = 24; D = "3.746.10";    E = "  +0.94%  " },
    @{ Row = 26; E = "  +1.27%  " },
    @{ Row = 27; E = "  +3.04%  " },
    @{ Row = 28; D = "9.30";        E = "  +0.94%  " },
    @{ Row = 29; E = "  -0.19%  " },
    @{ Row = 30; E = "  +0.00%  " },
    @{ Row = 31; D = "3.594.93";    E = "  +0.87%  " },
    @{ Row = 32; E = "  -1.82%  " },
    @{ Row = 33; D = "25.49";       E = "  +0.16%  " },
    @{ Row = 34; E = "  -2.28%  " },
    @{ Row = 35; D = "7.84";        E = "  -0.89%  " },
    @{ Row = 36; E = "  +0.00%  " },
    @{ Row = 37; E = "  -2.66%  " },
    @{ Row = 38; E = "  -0.23%  " },
    @{ Row = 39; D = "176.77";      E = "  +2.98%  " },
    @{ Row = 40; D = "0.0859";      E = "  +0.20%  " },
    @{ Row = 41; D = "5.23";        E = "  +0.19%  " },
    @{ Row = 42; D = "0.898";       E = "  +0.02%  " },
    @{ Row = 43; E = "  -2.47%  " },
    @{ Row = 44; D = "2.58";        E = "  +8.05%  " },
    @{ Row = 45; D = "1.00";        E = "  +0.12%  " },
    @{ Row = 46; D = "1.19";        E = "  -2.00%  " },
    @{ Row = 49; D = "7.18";        E = "  +0.48%  " },
    @{ Row = 50; D = "0.951" },
    @{ Row = 51; D = "2.415.95";    E = "  +4.66%  " }
)

foreach ($r in $rows) {
    if ($r.ContainsKey("D")) {
        $addr = "D" + $r.Row
        $ws.Range($addr).Value = "'" + $r.D
        $ws.Range($addr).Style = "Normal"
    }
    if ($r.ContainsKey("E")) {
        $addr = "E" + $r.Row
        $ws.Range($addr).Value = "'" + $r.E
        $ws.Range($addr).Style = "Normal"
    }
}

# Rows 47/48 swap places in the ranking: InjectiveProtocol drops to #48,
# EnergySwap rises to #47 -- update name, link, price and change together.
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'24.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +4.24%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'25.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.88%  "
$ws.Range("E48").Style = "Normal"
